$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is numeric-looking text (e.g. '209.42') need to be
# pre-formatted as Text so Excel stores them as strings (matching the source
# workbook's inlineStr cells) instead of coercing them into numbers.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D13", "D14", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D30", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '89.319.11'
$ws.Range("E2").Value = '  +0.67%  '

# Row 3
$ws.Range("D3").Value = '3.025.13'
$ws.Range("E3").Value = '  -3.88%  '

# Row 4
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.13%  '

# Row 5
$ws.Range("D5").Value = '209.42'
$ws.Range("E5").Value = '  -2.86%  '

# Row 6
$ws.Range("D6").Value = '610.93'
$ws.Range("E6").Value = '  -3.66%  '

# Row 7
$ws.Range("D7").Value = '0.358'
$ws.Range("E7").Value = '  -10.19%  '

# Row 8
$ws.Range("D8").Value = '0.877'
$ws.Range("E8").Value = '  +19.91%  '

# Row 9
$ws.Range("E9").Value = '  +0.10%  '

# Row 10
$ws.Range("D10").Value = '3.021.47'
$ws.Range("E10").Value = '  -3.98%  '

# Row 11
$ws.Range("E11").Value = '  +19.19%  '

# Row 12
$ws.Range("E12").Value = '  +3.80%  '

# Row 13
$ws.Range("D13").Value = '0.0000235'
$ws.Range("E13").Value = '  -6.79%  '

# Row 14
$ws.Range("D14").Value = '5.32'
$ws.Range("E14").Value = '  +0.41%  '

# Row 15
$ws.Range("D15").Value = '88.711.51'

# Row 16
$ws.Range("D16").Value = '31.67'
$ws.Range("E16").Value = '  -2.68%  '

# Row 17
$ws.Range("D17").Value = '3.592.93'
$ws.Range("E17").Value = '  -3.11%  '

# Row 18
$ws.Range("D18").Value = '3.026.79'
$ws.Range("E18").Value = '  -3.49%  '

# Row 19
$ws.Range("D19").Value = '3.33'
$ws.Range("E19").Value = '  -0.40%  '

# Row 20
$ws.Range("D20").Value = '0.0000209'
$ws.Range("E20").Value = '  -8.37%  '

# Row 21
$ws.Range("D21").Value = '13.27'
$ws.Range("E21").Value = '  -0.10%  '

# Row 22
$ws.Range("D22").Value = '421.68'
$ws.Range("E22").Value = '  -1.24%  '

# Row 23
$ws.Range("D23").Value = '4.96'
$ws.Range("E23").Value = '  +1.05%  '

# Row 24
$ws.Range("D24").Value = '8.08'
$ws.Range("E24").Value = '  -3.82%  '

# Row 25
$ws.Range("D25").Value = '5.33'
$ws.Range("E25").Value = '  -1.64%  '

# Row 26
$ws.Range("D26").Value = '82.88'
$ws.Range("E26").Value = '  +3.53%  '

# Row 27
$ws.Range("D27").Value = '11.58'
$ws.Range("E27").Value = '  +0.36%  '

# Row 28
$ws.Range("D28").Value = '3.195.61'
$ws.Range("E28").Value = '  -2.59%  '

# Row 29
$ws.Range("E29").Value = '  -0.03%  '

# Row 30
$ws.Range("D30").Value = '0.161'
$ws.Range("E30").Value = '  +1.64%  '

# Row 31
$ws.Range("E31").Value = '  +1.57%  '

# Row 32
$ws.Range("D32").Value = '8.11'
$ws.Range("E32").Value = '  -1.47%  '

# Row 33
$ws.Range("D33").Value = '499.20'
$ws.Range("E33").Value = '  -2.74%  '

# Row 34
$ws.Range("D34").Value = '3.57'
$ws.Range("E34").Value = '  -11.36%  '

# Row 35
$ws.Range("D35").Value = '6.57'
$ws.Range("E35").Value = '  -7.97%  '

# Row 36
$ws.Range("D36").Value = '22.61'
$ws.Range("E36").Value = '  +3.13%  '

# Row 37
$ws.Range("E37").Value = '  -2.81%  '

# Row 38
$ws.Range("D38").Value = '1.22'
$ws.Range("E38").Value = '  -6.50%  '

# Row 39
$ws.Range("B39").Value = 'WhiteBITCoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D39").Value = '22.21'
$ws.Range("E39").Value = '  -0.09%  '

# Row 40
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = '0.129'
$ws.Range("E40").Value = '  -4.94%  '

# Row 41
$ws.Range("E41").Value = '  -0.01%  '

# Row 42
$ws.Range("E42").Value = '  -0.02%  '

# Row 43
$ws.Range("D43").Value = '0.359'
$ws.Range("E43").Value = '  -1.54%  '

# Row 44
$ws.Range("D44").Value = '0.136'
$ws.Range("E44").Value = '  +7.67%  '

# Row 45
$ws.Range("D45").Value = '1.81'
$ws.Range("E45").Value = '  -3.35%  '

# Row 46
$ws.Range("D46").Value = '145.61'
$ws.Range("E46").Value = '  -0.20%  '

# Row 47
$ws.Range("D47").Value = '43.25'
$ws.Range("E47").Value = '  -1.03%  '

# Row 48
$ws.Range("D48").Value = '0.0669'
$ws.Range("E48").Value = '  +8.33%  '

# Row 49
$ws.Range("D49").Value = '4.02'
$ws.Range("E49").Value = '  +1.98%  '

# Row 50
$ws.Range("D50").Value = '1.20'
$ws.Range("E50").Value = '  +0.87%  '

# Row 51
$ws.Range("D51").Value = '157.94'
$ws.Range("E51").Value = '  -4.91%  '

# Restore default (General/Normal) formatting on those cells now that the text
# values are stored, so no residual cell-level style difference remains.
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).Style = "Normal"
}
